# Actualización automática 2025-06-02 13:21:56
# Adds a new "PRESUPUESTO" column (G) to the "VENTA MENSUAL" sheet,
# mirroring the existing monthly columns' formatting, filled with 0s
# for every data row and matching the bottom totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column width (matches <col width="17" customWidth="1" min="7" max="7"/>)
# Note: Excel's COM ColumnWidth is ~0.83 narrower than the stored OOXML
# <col width>, so subtract that offset to land on exactly 17 in the XML.
$ws.Columns.Item(7).ColumnWidth = 16.17

# Header cell G1 = "PRESUPUESTO", styled like the other header cells (A1:F1)
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows G2:G28 -> 0, styled like the existing currency columns (C2:C28)
$ws.Range("C2").Copy()
$ws.Range("G2:G28").PasteSpecial(-4122)
$ws.Range("G2:G28").Value = 0

# Totals row G29 -> 0, styled like the existing totals row (C29)
$ws.Range("C29").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 0
